# DaySale_2025-08-30_00-00.xlsx update:
#  - insert a new item row ("OMEZ 20MG 14 CAPS.") into the shortages table
#  - renumber the row that follows it (سرنجات 3 سم) from 3 to 4
#  - refresh the grand-total cell
#  - bump the "generated at" timestamp in the footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert a blank row above the existing row 9 ("سرنجات 3 سم").
#    Everything at/after row 9 (that row, the totals row, the footer row)
#    shifts down by one, carrying its formatting/merges with it.
$ws.Rows.Item(9).Insert()

# 2) Clone the per-column formatting of the item-row template (now sitting
#    one row down, at row 10) onto the freshly inserted row 9.
$ws.Range("A10:Q10").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# 3) Re-create the merged cells for the new row 9 (PasteSpecial of formats
#    does not carry merge state along with it).
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# 4) Fill in the new item's data.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "OMEZ 20MG 14 CAPS."

$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "0:1"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"

$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "56.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "28.0000"

$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "0:1"

# Forcing "@" text formats above clobbers the template's original number
# formats; restore those (display only - the cell contents already hold
# the right text) without disturbing the values just entered.
$ws.Range("H10:Q10").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# Row heights follow the sheet's row-position template, not the row's
# moving content, so set them explicitly to match.
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75

# 5) The row that used to be #3 in the list ("سرنجات 3 سم") is now #4.
$ws.Range("A10").Value = 4

# 6) Refresh the grand total (now on row 11 after the insert) and restore
#    its template row height.
$ws.Range("P11").Value = 150.75999999999999
$ws.Rows.Item(11).RowHeight = 25.5

# 7) Bump the generated-at timestamp shown in the footer (now row 12).
$ws.Range("A12").Value = "Saturday, 30 August, 2025 10:27 AM"
